$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the edited cells keep a Text format so values (VAT numbers, dates)
# are not auto-converted into numbers/dates by Excel.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("C2").NumberFormat = "@"
$ws.Range("F2").NumberFormat = "@"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("C3").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"

# Row 2 updates
$ws.Range("B2").Value = "Integrated Steel Company for Cable Tray System(ISP)"
$ws.Range("C2").Value = "31080000057900003"
$ws.Range("F2").Value = "2398.00"

# Row 3 updates
$ws.Range("B3").Value = "Al-Hawamah Trading Est."
$ws.Range("C3").Value = "3020709160003"
$ws.Range("D3").Value = "2022/08/07"
